$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tiny floating-point correction on the existing last row (A17) -
# this matches the re-saved serial date value for 2025-07-27 18:00:26.
$ws.Range("A17").Value = 45865.75030394676

# Append the new row 18 with the latest scheduled-task reading
# (2025-07-27 19:00:26 execution).
$ws.Range("A18").Value = 45865.79197570039
$ws.Range("A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B18").Value = 2025
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 14.78
$ws.Range("E18").Value = 87.33
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 8.43
$ws.Range("H18").Value = "E"
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "19:00:26"
